$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new daily data rows (44-50) ---

# Row 44
$ws.Range("A44").Value = 43368
$ws.Range("D44").Formula = "=(C44-B44)* 1440"
$ws.Range("E44").Formula = "=IF(C44>B44, (C44-B44)*1440, (B44-C44)*1440)"
$ws.Range("F44").Formula = "=ABS((C44-B44)*1440)"

# Row 45
$ws.Range("A45").Value = 43369
$ws.Range("D45").Formula = "=(C45-B45)* 1440"
$ws.Range("E45").Formula = "=IF(C45>B45, (C45-B45)*1440, (B45-C45)*1440)"
$ws.Range("F45").Formula = "=ABS((C45-B45)*1440)"

# Row 46
$ws.Range("A46").Value = 43370
$ws.Range("D46").Formula = "=(C46-B46)* 1440"
$ws.Range("E46").Formula = "=IF(C46>B46, (C46-B46)*1440, (B46-C46)*1440)"
$ws.Range("F46").Formula = "=ABS((C46-B46)*1440)"

# Row 47
$ws.Range("A47").Value = 43371
$ws.Range("B47").Value = 0
$ws.Range("C47").Value = 0
$ws.Range("D47").Formula = "=(C47-B47)* 1440"
$ws.Range("E47").Formula = "=IF(C47>B47, (C47-B47)*1440, (B47-C47)*1440)"
$ws.Range("F47").Formula = "=ABS((C47-B47)*1440)"

# Row 48
$ws.Range("A48").Value = 43372
$ws.Range("B48").Value = 0
$ws.Range("C48").Value = 0
$ws.Range("D48").Formula = "=(C48-B48)* 1440"
$ws.Range("E48").Formula = "=IF(C48>B48, (C48-B48)*1440, (B48-C48)*1440)"
$ws.Range("F48").Formula = "=ABS((C48-B48)*1440)"

# Row 49
$ws.Range("A49").Value = 43373
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("D49").Formula = "=(C49-B49)* 1440"
$ws.Range("E49").Formula = "=IF(C49>B49, (C49-B49)*1440, (B49-C49)*1440)"
$ws.Range("F49").Formula = "=ABS((C49-B49)*1440)"

# Row 50 (formula-only row, no date/time entered yet)
$ws.Range("D50").Formula = "=(C50-B50)* 1440"
$ws.Range("E50").Formula = "=IF(C50>B50, (C50-B50)*1440, (B50-C50)*1440)"
$ws.Range("F50").Formula = "=ABS((C50-B50)*1440)"

# Match the direct (non-column-inherited) cell formatting used on the
# Second Duration / Absolute Value calculated columns for these new rows.
$ws.Range("E44:F50").NumberFormat = "General"

# Grow the table to cover the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F50"))

# --- Update the view state to match where the user ended up working ---
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("B45").Select()
